$d = $word.ActiveDocument

$replacements = @(
    @("{{dissertation_topic}}", "[%dissertation_topic%]"),
    @("{{student_specialty}}", "[%student_specialty%]"),
    @("{{student_full_name}}", "[%student_full_name%]"),
    @("{{student_program}}", "[%student_program%]"),
    @("{{student_phone}}", "[%student_phone%]"),
    @("{{student_email}}", "[%student_email%]"),
    @("{{day}}", "[%day%]"),
    @("{{month}}", "[%month%]"),
    @("{{year}}", "[%year%]")
)

foreach ($p in $d.Paragraphs) {
    $rng = $p.Range
    $original = $rng.Text

    # Paragraph ranges report their text including the trailing paragraph
    # mark (CR). Strip it off before editing so re-assigning .Text does not
    # insert an extra paragraph break.
    $hasMark = $original.Length -gt 0 -and $original.EndsWith([char]13)
    $core = $original
    if ($hasMark) {
        $core = $original.Substring(0, $original.Length - 1)
    }

    $updated = $core
    foreach ($pair in $replacements) {
        $updated = $updated.Replace($pair[0], $pair[1])
    }

    if ($updated -ne $core) {
        $rng.Text = $updated
    }
}
